$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq "6-6-2013-14") {
        # Assigning the literal text "2014-06-06" directly would be auto-
        # recognized as a date by Excel's input heuristics (it would turn
        # into a date serial number with a date number format applied).
        # Using a formula whose result is a quoted text string keeps the
        # value as plain text, then converting that formula to a static
        # value via copy / paste-special avoids touching the cell's style.
        $cell.Formula = '="2014-06-06"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$excel.CutCopyMode = 0
